# subsid.xlsx fix: update headers, fill "no data" placeholders, and
# renumber rows 18-24 (sequence continues 16,17,18,19,20,21,22 instead of
# restarting at 18).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row fixes ---
$ws.Range("C1").Value = "Муниципалитет"
$ws.Range("D1").Value = "Населенный пункт"

# --- Row 7: fill empty E:H cells with placeholder text ---
$ws.Range("E7:H7").Value = "нет данных"

# --- Rows 18-24: renumber column A and fill empty E:H cells ---
$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(20, 1).Value = 18
$ws.Cells.Item(21, 1).Value = 19
$ws.Cells.Item(22, 1).Value = 20
$ws.Cells.Item(23, 1).Value = 21
$ws.Cells.Item(24, 1).Value = 22

$ws.Range("E18:H24").Value = "нет данных"
